$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the underlying byte/half-word counts (commit: "store bytes instead of half word")
$ws.Range("B14").Value = 361
$ws.Range("B15").Value = 8

# Update the view state: scroll the window so column A is visible again, and
# move the active selection to L13 (matches the saved sheetView/selection).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("L13").Select()
